$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sir Tomato sprite update - add "Lava puddle" block type entry
$ws.Range("A14").Value = "Lava puddle"
$ws.Range("B14").Value = "l"

# Update the active cell selection
$ws.Range("A15").Select()
